$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data changes: rows 10-12, column C (Execute) and D (Invocation Count) ---
$ws.Range("C10").Value = "Yes"
$ws.Range("D10").Value = 2

$ws.Range("C11").Value = "Yes"
$ws.Range("D11").Value = 2

$ws.Range("D12").Value = 2

# --- Style change: A6 gets a fresh default-font style (new font/xf entry) ---
$ws.Range("A6").Font.Name = "Calibri"

# --- View changes: scroll sheet back to top-left (remove topLeftCell="A7") and move selection to F9 ---
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("F9").Select()
